# Add 2022-Q4 data: new detail sheet + updated summary sheet ("总计").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" (summary) sheet: insert a new row for 2022-Q4 above the existing
#    2022-Q3 row, shifting everything else down.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# The inserted row copies the header row's formatting (bold/border/center) -
# clear that from the 3 plain data columns so they look like the other data
# rows, then restore the bold index-style on column A only.
$summary.Range("B2:D2").ClearFormats()

$a2 = $summary.Range("A2")
$a2.Value = 0
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 19
$summary.Range("D2").Value = 4.84

# ---------------------------------------------------------------------------
# 2. New "2022-Q4" worksheet with per-fund detail, inserted right after
#    "总计" (i.e. before the existing "2022-Q3" tab).
# ---------------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q4Sheet = $wb.Worksheets.Add($q3Sheet)
$q4Sheet.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 0; $col -lt $headers.Length; $col++) {
    $cell = $q4Sheet.Cells.Item(1, $col + 2)
    $cell.Value = $headers[$col]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$rows = @(
    @(0, "000603", "易方达创新驱动灵活配置混合", "47.69", "93.25", "3.85", "1.8361", 10),
    @(1, "003961", "易方达瑞程灵活配置混合A", "13.38", "93.64", "4.37", "0.5847", 5),
    @(2, "003962", "易方达瑞程灵活配置混合C", "6.76", "93.64", "4.37", "0.2954", 5),
    @(3, "014271", "大成北交所两年定开混合A", "3.24", "68.93", "7.42", "0.2404", 3),
    @(4, "014273", "广发北交所精选两年定开混合A", "3.23", "83.79", "7.42", "0.2397", 2),
    @(5, "014279", "汇添富北交所创新精选两年定开混合A", "3.06", "94.24", "7.32", "0.2240", 1),
    @(6, "580008", "东吴新产业精选股票A", "4.38", "91.75", "4.98", "0.2181", 3),
    @(7, "011470", "东吴新产业精选混合C", "4.38", "91.75", "4.98", "0.2181", 3),
    @(8, "014294", "南方北交所精选两年定开混合", "4.05", "90.08", "5.34", "0.2163", 5),
    @(9, "014275", "易方达北交所精选两年定开混合A", "3.47", "70.57", "6.10", "0.2117", 1),
    @(10, "014283", "华夏北交所创新中小企业精选两年定开混合", "3.27", "90.95", "6.07", "0.1985", 5),
    @(11, "014269", "嘉实北交所精选两年定期混合A", "2.65", "94.48", "3.27", "0.0867", 10),
    @(12, "014274", "广发北交所精选两年定开混合C", "0.81", "83.79", "7.42", "0.0601", 2),
    @(13, "014272", "大成北交所两年定开混合C", "0.77", "68.93", "7.42", "0.0571", 3),
    @(14, "014276", "易方达北交所精选两年定开混合C", "0.90", "70.57", "6.10", "0.0549", 1),
    @(15, "016307", "景顺长城北交所精选两年定开混合A", "1.83", "43.56", "2.16", "0.0395", 8),
    @(16, "014280", "汇添富北交所创新精选两年定开混合C", "0.48", "94.24", "7.32", "0.0351", 1),
    @(17, "014270", "嘉实北交所精选两年定期混合C", "0.52", "94.48", "3.27", "0.0170", 10),
    @(18, "016308", "景顺长城北交所精选两年定开混合C", "0.27", "43.56", "2.16", "0.0058", 8)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $rowNum = $i + 2

    $cellA = $q4Sheet.Cells.Item($rowNum, 1)
    $cellA.Value = $r[0]
    $cellA.Font.Bold = $true
    $cellA.HorizontalAlignment = -4108
    $cellA.VerticalAlignment = -4160
    $cellA.Borders.LineStyle = 1

    $q4Sheet.Cells.Item($rowNum, 2).Value = "'" + $r[1]
    $q4Sheet.Cells.Item($rowNum, 3).Value = "'" + $r[2]
    $q4Sheet.Cells.Item($rowNum, 4).Value = "'" + $r[3]
    $q4Sheet.Cells.Item($rowNum, 5).Value = "'" + $r[4]
    $q4Sheet.Cells.Item($rowNum, 6).Value = "'" + $r[5]
    $q4Sheet.Cells.Item($rowNum, 7).Value = "'" + $r[6]
    $q4Sheet.Cells.Item($rowNum, 8).Value = $r[7]
}
